# Task #966: Add CostTypeExportWizard
# Rename the InterfaceType* sheets into the new Cost* naming scheme,
# update their titles/headers/sample data, and drop the no-longer-needed
# "Interfaces" sheet entirely.

$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsCostTypes = $wb.Worksheets.Item("InterfaceTypes")
$wsCostTypes.Name = "CostTypes"

$wsCostEquip = $wb.Worksheets.Item("InterfaceEnds")
$wsCostEquip.Name = "CostEquipments"

# --- Drop the Interfaces sheet -------------------------------------------
$wsInterfaces = $wb.Worksheets.Item("Interfaces")
[void]$wsInterfaces.Delete()

# --- Update "CostTypes" sheet content ------------------------------------
$wsCostTypes.Range("A2").Value = "VirSat IO Sheet for Cost Types"
$wsCostTypes.Range("C4").Value = "Cost Type Name"
$wsCostTypes.Range("C5").Value = "BIII"
$wsCostTypes.Range("C6").Value = "BAAA"
$wsCostTypes.Range("C7").Value = "BUUU"
[void]$wsCostTypes.Range("C4").Select()

# --- Update "CostEquipments" sheet content -------------------------------
$wsCostEquip.Range("A2").Value = "VirSat IO Sheet for CostEquipments"
$wsCostEquip.Range("C4").Value = "Cost Name"
$wsCostEquip.Range("D4").Value = "Cost Type"
[void]$wsCostEquip.Range("D4").Select()
